$d = $word.ActiveDocument

# Locate the Subtitle paragraph ("Punishment and self-governance among men
# serving life sentences for murder") so we can insert a new "Author"
# paragraph ("Ben Jarman") directly after it.
$subtitleText = "Punishment and self-governance among men serving life sentences for murder"

$targetPara = $null
foreach ($para in $d.Paragraphs) {
    $t = $para.Range.Text
    $t = $t.TrimEnd([char]13, [char]7)
    if ($t -eq $subtitleText) {
        $targetPara = $para
        break
    }
}

if ($targetPara -eq $null) {
    throw "Could not find the Subtitle paragraph to anchor the new Author paragraph"
}

# The insertion point must sit just before the paragraph mark of the
# Subtitle paragraph (i.e. End - 1), not exactly on the boundary shared
# with the next paragraph's Start (End), otherwise the inserted content
# gets folded into the following paragraph instead of forming its own.
$insertPos = $targetPara.Range.End - 1
$insertRange = $d.Range($insertPos, $insertPos)

$xml = @'
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="Author"/>
            </w:pPr>
            <w:r>
              <w:t xml:space="preserve">Ben Jarman</w:t>
            </w:r>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
'@

$insertRange.InsertXML($xml)
